# The sheet originally has:
#   A: segment name (text)   B..I: 8 numeric stat columns   J: totalMean   K: totalStd
# The target layout inserts a new leading "segments" column holding the
# 0-based segment index, pushing the segment-name column (and everything
# after it) one column to the right:
#   A: segment index (number) B: segment name (text)  C..J: stats  K: totalMean  L: totalStd

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column before column B - this shifts the old B..K data
# (and their styles) right into C..L, and (because it copies the format of
# the column immediately to its left, column A) gives the new column B the
# same style as column A for the data rows.
$ws.Columns("B:B").Insert()

# Header for the freshly inserted column.
$ws.Range("B1").Value = "segments"

# 0-based segment index (new column A) and the segment name that used to
# live in column A (now column B), for each of the 19 data rows.
$indices = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18)
$names = @(
    "background","back_bumper","back_glass","back_left_door","back_left_light",
    "back_right_door","back_right_light","front_bumper","front_glass","front_left_door",
    "front_left_light","front_right_door","front_right_light","hood","left_mirror",
    "right_mirror","tailgate","trunk","wheel"
)

for ($i = 0; $i -lt $indices.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $indices[$i]
    $ws.Cells.Item($r, 2).Value = $names[$i]
}

# The new "segments" name column is unstyled in the data rows (only the
# index column keeps the bordered/centered/bold style), so strip the style
# that Insert() copied into it from column A.
$ws.Range("B2:B20").ClearFormats()

# ...but the header cell B1 must match the other (styled) header cells, so
# copy formatting only (not value) from a neighbouring header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
